# Revert the "Invoice" sample/demo content back to a blank template:
#  - remove the placed-in-cell logo image that was occupying A1 (shows as
#    a #VALUE! rich-value error in the trial/unlicensed render)
#  - clear the sample line-item rows (Qty/Unit Price/Total, C15:E27) that
#    were only there to demo the invoice; the Subtotal/Tax/Total formulas
#    in row 29-31 are left in place and simply recompute to 0
#  - reset the saved selection / scroll position back to the top of the
#    sheet instead of the scrolled-down G20 selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")
$ws.Activate()

# Remove the embedded/placed image value that lived in A1.
$ws.Range("A1").ClearContents()

# Clear the demo line items (quantity, unit price, computed total) for
# every row of the items table - formulas included, so the cells go
# fully blank rather than recalculating to 0.
$ws.Range("C15:E27").ClearContents()

# Recalculate so the Subtotal / Tax / Total formulas (which still
# reference E15:E27) pick up the now-empty inputs.
$excel.Calculate()

# Reset the view: scroll back to the top-left and select A2 instead of
# the previous scrolled-down G20 selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()

# Best-effort: put the window back into Normal view (some engines may
# not persist this, but it mirrors the source edit's removal of
# view="pageLayout"/topLeftCell from the sheet view).
$excel.ActiveWindow.View = 1

# Best-effort: restore the workbook window's on-screen X position.
$wb.Windows.Item(1).Left = 0
